# V0.8 Basic domain validation
# Adds extra domain columns (F:K) to the header/test rows and marks the
# cells that contain a valid, well-formed URL as Excel hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cells F1:K1
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "github.com"
$ws.Range("G1").Value = ".wrong.url"
$ws.Range("H1").Value = "www.test.com"
$ws.Range("I1").Value = "bbc.com/subsite"
$ws.Range("J1").Value = "https://wwww.google.com"
$ws.Range("K1").Value = "https://google.com"

# Give the new plain-text header cells the same look as the existing
# header cells (B1:E1): bold, centered, rotated 180, wrapped, text format.
$plainHeaders = $ws.Range("F1:I1")
$plainHeaders.NumberFormat = "@"
$plainHeaders.Font.Bold = $true
$plainHeaders.HorizontalAlignment = -4108
$plainHeaders.Orientation = -4170
$plainHeaders.WrapText = $true

# H1, J1 and K1 contain well formed URLs -> turn them into real hyperlinks.
# Excel's own "Hyperlink" cell style (underline, theme color) is applied
# automatically by Hyperlinks.Add.
$ws.Hyperlinks.Add($ws.Range("H1"), "http://www.test.com", $null, $null, $null)
$ws.Hyperlinks.Add($ws.Range("J1"), "https://wwww.google.com", $null, $null, $null)
$ws.Hyperlinks.Add($ws.Range("K1"), "https://google.com", $null, $null, $null)

$hyperlinkHeaders = $ws.Range("H1,J1,K1")
$hyperlinkHeaders.HorizontalAlignment = -4108
$hyperlinkHeaders.Orientation = -4170
$hyperlinkHeaders.WrapText = $true

# ---------------------------------------------------------------------
# 2. New "x" marker cells on the sample rows
# ---------------------------------------------------------------------
$row2 = $ws.Range("F2:K2")
$row2.Value = "x"
$row2.NumberFormat = "@"
$row2.HorizontalAlignment = -4108

$ws.Range("F5").Value = "x"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").HorizontalAlignment = -4108

$row6 = $ws.Range("G6:I6")
$row6.Value = "x"
$row6.NumberFormat = "@"
$row6.HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3. Window / view changes
# ---------------------------------------------------------------------
$wb.Windows.Item(1).Left = 6135
$excel.ActiveWindow.Zoom = 175

$ws.Range("K4").Select() | Out-Null
